# Applies the "final tweaks, minified CSS, updated report" edit:
#  - Rewrites the four "Before/After corrections" caption paragraphs
#    (Gigs page x2, Recommendations page x2) with new wording/run layout.
#  - Adds a <w:lastRenderedPageBreak/> marker to the "AFTER" Recommendations
#    caption.
#  - Removes one extra blank paragraph before each Recommendations caption.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: Paragraph.Range.Text includes the trailing paragraph-mark
# character (CR, chr 13) - strip it so text comparisons are exact.
# ---------------------------------------------------------------------
function Get-ParaText($para) {
    $raw = $para.Range.Text
    $clean = $raw.TrimEnd([char]13)
    return $clean
}

# ---------------------------------------------------------------------
# Helper: find the 1-based Paragraphs index of the (first) paragraph whose
# text equals $text exactly.
# ---------------------------------------------------------------------
function Find-ParagraphIndex($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $ptext = Get-ParaText $para
        if ($ptext -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Helper: replace the contents of paragraph $index with the run structure
# described by the inner-body XML fragment $innerXml (one or more <w:r>
# elements, no surrounding <w:p>). Preserves the paragraph mark / paraId.
# ---------------------------------------------------------------------
function Set-ParagraphRuns($doc, $index, $innerXml) {
    $para = $doc.Paragraphs.Item($index)
    $full = $para.Range
    $bodyStart = $full.Start
    $bodyEnd = $full.End - 1
    # Delete everything except the trailing paragraph-mark character.
    if ($bodyEnd -gt $bodyStart) {
        $body = $doc.Range($bodyStart, $bodyEnd)
        $body.Delete()
    }

    $para2 = $doc.Paragraphs.Item($index)
    $insertStart = $para2.Range.Start
    $insertAt = $doc.Range($insertStart, $insertStart)

    $pkgHead = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p>"
    $pkgTail = "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $pkg = $pkgHead + $innerXml + $pkgTail

    $insertAt.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) "Before corrections for Gigs page"
#    -> "Gigs page HTML file " | "B" | "EFORE" | " corrections" | ":" | " "
# ---------------------------------------------------------------------
$target1 = "Before corrections for Gigs page"
$idx1 = Find-ParagraphIndex $d $target1
$xml1 = "<w:r><w:t xml:space='preserve'>Gigs page HTML file </w:t></w:r><w:r><w:t>B</w:t></w:r><w:r><w:t>EFORE</w:t></w:r><w:r><w:t xml:space='preserve'> corrections</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space='preserve'> </w:t></w:r>"
Set-ParagraphRuns $d $idx1 $xml1

# ---------------------------------------------------------------------
# 2) "After corrections of Gigs page:"
#    -> single run "Gigs page HTML file AFTER corrections: "
# ---------------------------------------------------------------------
$target2 = "After corrections of Gigs page:"
$idx2 = Find-ParagraphIndex $d $target2
$xml2 = "<w:r><w:t xml:space='preserve'>Gigs page HTML file AFTER corrections: </w:t></w:r>"
Set-ParagraphRuns $d $idx2 $xml2

# ---------------------------------------------------------------------
# 3) "Before corrections for Recommendations page:"
#    -> "Recommendations" | " page HTML file BEFORE corrections: "
#    Also: remove one of the five blank paragraphs immediately preceding it.
# ---------------------------------------------------------------------
$target3 = "Before corrections for Recommendations page:"
$idx3pre = Find-ParagraphIndex $d $target3
$blankIdx3 = $idx3pre - 1
$blank3 = $d.Paragraphs.Item($blankIdx3)
$blank3text = Get-ParaText $blank3
if ($blank3text -eq "") {
    $blank3.Range.Delete()
}

$idx3 = Find-ParagraphIndex $d $target3
$xml3 = "<w:r><w:t>Recommendations</w:t></w:r><w:r><w:t xml:space='preserve'> page HTML file BEFORE corrections: </w:t></w:r>"
Set-ParagraphRuns $d $idx3 $xml3

# ---------------------------------------------------------------------
# 4) "After corrections of Recommendations page"
#    -> <w:lastRenderedPageBreak/>"Recommendations page HTML file " | "AFTER" | " corrections: "
#    Also: remove one of the five blank paragraphs immediately preceding it.
# ---------------------------------------------------------------------
$target4 = "After corrections of Recommendations page"
$idx4pre = Find-ParagraphIndex $d $target4
$blankIdx4 = $idx4pre - 1
$blank4 = $d.Paragraphs.Item($blankIdx4)
$blank4text = Get-ParaText $blank4
if ($blank4text -eq "") {
    $blank4.Range.Delete()
}

$idx4 = Find-ParagraphIndex $d $target4
$xml4 = "<w:r><w:lastRenderedPageBreak/><w:t xml:space='preserve'>Recommendations page HTML file </w:t></w:r><w:r><w:t>AFTER</w:t></w:r><w:r><w:t xml:space='preserve'> corrections: </w:t></w:r>"
Set-ParagraphRuns $d $idx4 $xml4

Write-Output "done"
